# Generate Report for Handoff
# Refresh the localization-status report: the source file
# 17295b83-c8e9-45e2-9fe8-940e2078275f.md was replaced by
# 726fd1a3-c1f3-40c1-87ca-a25100b79da0.md, two new image dependencies
# (.png) were discovered and handed off alongside it, and the handoff
# timestamps/handoff-file hashes were refreshed.

$wb = $excel.ActiveWorkbook

$mdName   = "726fd1a3-c1f3-40c1-87ca-a25100b79da0.md"
$png1Name = "d35dacc8-e8de-4a5d-aa10-287ee8165bb6.png"
$png2Name = "e57ee89c-f0e2-4b4e-90f9-0bf3bab42277.png"

$zhXlfName  = "726fd1a3-c1f3-40c1-87ca-a25100b79da0.cd64925ca30223789e10d3abe7d117930103f97e.zh-cn.xlf"
$deXlfName  = "726fd1a3-c1f3-40c1-87ca-a25100b79da0.cd64925ca30223789e10d3abe7d117930103f97e.de-de.xlf"
$png1Target = "f142bd1157416b636cf837b39c6655bb92cdf10e.png"
$png2Target = "1ee8bf0e08142d7ef08429ffc59e718a834840a9.png"

$handoffDate    = "2016-03-24 03:10:32"
$zhHandoffDate  = "2016-03-24 03:10:28"
$epoch          = "0001-01-01 00:00:00"
$dependencyFrom = "e2e\726fd1a3-c1f3-40c1-87ca-a25100b79da0.md"

$mdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/726fd1a3c1f340c187caa25100b79da0cbf2eaf/e2e/$mdName"
$png1Url = "https://github.com/OpenLocalizationTest/oltest/blob/726fd1a3c1f340c187caa25100b79da0cbf2eaf/e2e/$png1Name"
$png2Url = "https://github.com/OpenLocalizationTest/oltest/blob/726fd1a3c1f340c187caa25100b79da0cbf2eaf/e2e/$png2Name"

$zhXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cd64925ca30223789e10d3abe7d117930103f97e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlfName"
$deXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cd64925ca30223789e10d3abe7d117930103f97e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlfName"
$png1ZhUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cd64925ca30223789e10d3abe7d117930103f97e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$png1Target"
$png1DeUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cd64925ca30223789e10d3abe7d117930103f97e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$png1Target"
$png2ZhUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cd64925ca30223789e10d3abe7d117930103f97e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$png2Target"
$png2DeUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cd64925ca30223789e10d3abe7d117930103f97e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$png2Target"

function Set-HyperlinkCell($ws, $cell, $text, $url) {
    $cell.Value = $text
    $ws.Hyperlinks.Add($cell, $url, "", "", $text) | Out-Null
}

# ---------------------------------------------------------------------
# Sheet "Overview" -- refresh row 2, then append rows 3 and 4
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-HyperlinkCell $wsOverview $wsOverview.Cells.Item(2, 1) $mdName $mdUrl
$wsOverview.Cells.Item(2, 2).Value = "Ready for handoff"
$wsOverview.Cells.Item(2, 3).Value = "Ready for handoff"
$wsOverview.Cells.Item(2, 4).Value = $handoffDate
$wsOverview.Cells.Item(2, 4).NumberFormat = "yyyy-mm-dd HH:mm:ss"

Set-HyperlinkCell $wsOverview $wsOverview.Cells.Item(3, 1) $png1Name $png1Url
$wsOverview.Cells.Item(3, 2).Value = "Ready for handoff"
$wsOverview.Cells.Item(3, 3).Value = "Ready for handoff"
$wsOverview.Cells.Item(3, 4).Value = $handoffDate
$wsOverview.Cells.Item(3, 4).NumberFormat = "yyyy-mm-dd HH:mm:ss"

Set-HyperlinkCell $wsOverview $wsOverview.Cells.Item(4, 1) $png2Name $png2Url
$wsOverview.Cells.Item(4, 2).Value = "Ready for handoff"
$wsOverview.Cells.Item(4, 3).Value = "Ready for handoff"
$wsOverview.Cells.Item(4, 4).Value = $handoffDate
$wsOverview.Cells.Item(4, 4).NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "zh-cn" -- refresh row 2, then append rows 3 and 4
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

Set-HyperlinkCell $wsZh $wsZh.Cells.Item(2, 1) $mdName $mdUrl
$wsZh.Cells.Item(2, 2).Value = ".md"
$wsZh.Cells.Item(2, 3).Value = "Ready for handoff"
Set-HyperlinkCell $wsZh $wsZh.Cells.Item(2, 4) $zhXlfName $zhXlfUrl
$wsZh.Cells.Item(2, 5).Value = $zhHandoffDate
$wsZh.Cells.Item(2, 8).Value = $epoch
$wsZh.Cells.Item(2, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(2, 10).Value = "Include"

Set-HyperlinkCell $wsZh $wsZh.Cells.Item(3, 1) $png1Name $png1Url
$wsZh.Cells.Item(3, 2).Value = ".png"
$wsZh.Cells.Item(3, 3).Value = "Ready for handoff"
Set-HyperlinkCell $wsZh $wsZh.Cells.Item(3, 4) $png1Target $png1ZhUrl
$wsZh.Cells.Item(3, 5).Value = $zhHandoffDate
$wsZh.Cells.Item(3, 8).Value = $epoch
$wsZh.Cells.Item(3, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(3, 10).Value = "IsDependency"
$wsZh.Cells.Item(3, 11).Value = $dependencyFrom

Set-HyperlinkCell $wsZh $wsZh.Cells.Item(4, 1) $png2Name $png2Url
$wsZh.Cells.Item(4, 2).Value = ".png"
$wsZh.Cells.Item(4, 3).Value = "Ready for handoff"
Set-HyperlinkCell $wsZh $wsZh.Cells.Item(4, 4) $png2Target $png2ZhUrl
$wsZh.Cells.Item(4, 5).Value = $zhHandoffDate
$wsZh.Cells.Item(4, 8).Value = $epoch
$wsZh.Cells.Item(4, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(4, 10).Value = "IsDependency"
$wsZh.Cells.Item(4, 11).Value = $dependencyFrom

# ---------------------------------------------------------------------
# Sheet "de-de" -- refresh row 2, then append rows 3 and 4
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

Set-HyperlinkCell $wsDe $wsDe.Cells.Item(2, 1) $mdName $mdUrl
$wsDe.Cells.Item(2, 2).Value = ".md"
$wsDe.Cells.Item(2, 3).Value = "Ready for handoff"
Set-HyperlinkCell $wsDe $wsDe.Cells.Item(2, 4) $deXlfName $deXlfUrl
$wsDe.Cells.Item(2, 5).Value = $handoffDate
$wsDe.Cells.Item(2, 5).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(2, 8).Value = $epoch
$wsDe.Cells.Item(2, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(2, 10).Value = "Include"

Set-HyperlinkCell $wsDe $wsDe.Cells.Item(3, 1) $png1Name $png1Url
$wsDe.Cells.Item(3, 2).Value = ".png"
$wsDe.Cells.Item(3, 3).Value = "Ready for handoff"
Set-HyperlinkCell $wsDe $wsDe.Cells.Item(3, 4) $png1Target $png1DeUrl
$wsDe.Cells.Item(3, 5).Value = $handoffDate
$wsDe.Cells.Item(3, 5).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(3, 8).Value = $epoch
$wsDe.Cells.Item(3, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(3, 10).Value = "IsDependency"
$wsDe.Cells.Item(3, 11).Value = $dependencyFrom

Set-HyperlinkCell $wsDe $wsDe.Cells.Item(4, 1) $png2Name $png2Url
$wsDe.Cells.Item(4, 2).Value = ".png"
$wsDe.Cells.Item(4, 3).Value = "Ready for handoff"
Set-HyperlinkCell $wsDe $wsDe.Cells.Item(4, 4) $png2Target $png2DeUrl
$wsDe.Cells.Item(4, 5).Value = $handoffDate
$wsDe.Cells.Item(4, 5).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(4, 8).Value = $epoch
$wsDe.Cells.Item(4, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(4, 10).Value = "IsDependency"
$wsDe.Cells.Item(4, 11).Value = $dependencyFrom
